# Applies the "Update addBudgetSequenceDiagram and ui class diagram" edit:
#  - Adds two new dashed/dotted Freeform connector shapes (copies of the
#    existing "Freeform 117" shape) pointing toward the ReportWindow area.
#  - Moves the "ReportWindow" rectangle (currently sandwiched between a
#    Freeform shape and an Elbow Connector) so that it is rendered on top
#    of (after, in z-order) the two new freeform connectors and the
#    ReportData box/connector that already sit at the end of the shape
#    stack.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> point helper (PowerPoint COM shape geometry is expressed in points;
# the OOXML stores EMUs, where 1 pt == 12700 EMU).
function EmuToPt([double]$emu) {
    return $emu / 12700.0
}

# --- Locate the template shape: the existing "Freeform 117" shape that
#     sits just above the ReportWindow rectangle (off 4116309,4167909 EMU /
#     ext 2642195,101600 EMU -> 324.1188,328.1818 pt). We match on name +
#     position so the script is resilient to shape ordering.
$templateFreeform = $null
$reportWindow = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.Name -eq "Freeform 117" -and `
        [math]::Abs($sh.Left - (EmuToPt 4116309)) -lt 0.01 -and `
        [math]::Abs($sh.Top - (EmuToPt 4167909)) -lt 0.01) {
        $templateFreeform = $sh
    }
    if ($sh.Name -eq "Rectangle 11" -and $sh.HasTextFrame -and `
        $sh.TextFrame.TextRange.Text -eq "ReportWindow") {
        $reportWindow = $sh
    }
}

if ($templateFreeform -eq $null) { throw "template Freeform 117 shape not found" }
if ($reportWindow -eq $null) { throw "ReportWindow shape not found" }

# --- New Freeform shape #1 (dashed connector ending near the ReportWindow
#     row, slightly below/right of the template).
$newFreeform1 = $templateFreeform.Duplicate()
$newFreeform1.Left = EmuToPt 4093862
$newFreeform1.Top = EmuToPt 5974080
$newFreeform1.Width = EmuToPt 2642195
$newFreeform1.Height = EmuToPt 101600

# --- New Freeform shape #2 (dashed connector, wider/shorter, a bit higher).
$newFreeform2 = $templateFreeform.Duplicate()
$newFreeform2.Left = EmuToPt 3492353
$newFreeform2.Top = EmuToPt 5539376
$newFreeform2.Width = EmuToPt 3243163
$newFreeform2.Height = EmuToPt 46175

# --- Move ReportWindow to the front so it renders after (on top of) the
#     two new freeform connectors that were just appended.
$reportWindow.ZOrder(0)   # msoBringToFront
